$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated numeric rows (theta, lambda, proportion drinking)
$ws.Range("C2").Value = 9.24
$ws.Range("D2").Value = 8.85
$ws.Range("E2").Value = 0.180291

$ws.Range("C4").Value = 10.28
$ws.Range("D4").Value = 10.6
$ws.Range("E4").Value = 0.138936

$ws.Range("C6").Value = 10.8
$ws.Range("D6").Value = 11.53
$ws.Range("E6").Value = 0.10858

$ws.Range("C8").Value = 9.6
$ws.Range("D8").Value = 11.32
$ws.Range("E8").Value = 0.106225

$ws.Range("C10").Value = 10.74
$ws.Range("D10").Value = 13.23
$ws.Range("E10").Value = 0.091207

$ws.Range("C12").Value = 11.63
$ws.Range("D12").Value = 13.11
$ws.Range("E12").Value = 0.090977

$ws.Range("C14").Value = 8.619999999999999
$ws.Range("D14").Value = 8.58
$ws.Range("E14").Value = 0.101183

# Rows that previously held the single "(nan)" text now hold distinct
# pairwise-estimate strings, one per cell. Shared strings are appended in
# column-major order (all of column C, then D, then E) to match the
# target workbook's shared-string table ordering.
$ws.Range("C3").Value = "(0.65)"
$ws.Range("C5").Value = "(0.69)"
$ws.Range("C7").Value = "(0.58)"
$ws.Range("C9").Value = "(0.3)"
$ws.Range("C11").Value = "(0.31)"
$ws.Range("C13").Value = "(1.01)"
$ws.Range("C15").Value = "(1.34)"

$ws.Range("D3").Value = "(0.24)"
$ws.Range("D5").Value = "(0.73)"
$ws.Range("D7").Value = "(0.18)"
$ws.Range("D9").Value = "(0.51)"
$ws.Range("D11").Value = "(0.26)"
$ws.Range("D13").Value = "(1.12)"
$ws.Range("D15").Value = "(1.02)"

$ws.Range("E3").Value = "(0.00502)"
$ws.Range("E5").Value = "(0.00963)"
$ws.Range("E7").Value = "(0.00581)"
$ws.Range("E9").Value = "(0.00552)"
$ws.Range("E11").Value = "(0.00342)"
$ws.Range("E13").Value = "(0.00503)"
$ws.Range("E15").Value = "(0.00851)"
